$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2063233333333333
$ws.Range("H2").Value = 0.61897
$ws.Range("M2").Value = 6.391557333333332
$ws.Range("N2").Value = 19.174672
$ws.Range("O2").Value = 0.1156448793857254
$ws.Range("P2").Value = 0.1156448793857254
$ws.Range("Q2").Value = 1.318727414204444
$ws.Range("R2").Value = 11.86854672784
$ws.Range("S2").Value = 0.1156448793857254
$ws.Range("T2").Value = 0.1156448793857254

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2063233333333333
$ws.Range("H3").Value = 0.61897
$ws.Range("O3").Value = 0.2610362896883882
$ws.Range("P3").Value = 0.2610362896883882
$ws.Range("Q3").Value = 2.976661942515555
$ws.Range("R3").Value = 26.78995748264
$ws.Range("S3").Value = 0.2610362896883882
$ws.Range("T3").Value = 0.2610362896883882

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2063233333333333
$ws.Range("H4").Value = 0.61897
$ws.Range("M4").Value = 32.348972
$ws.Range("N4").Value = 97.046916
$ws.Range("O4").Value = 0.5853022620452971
$ws.Range("P4").Value = 0.5853022620452972
$ws.Range("Q4").Value = 6.674347732946666
$ws.Range("R4").Value = 60.06912959652
$ws.Range("S4").Value = 0.5853022620452971
$ws.Range("T4").Value = 0.5853022620452972

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2063233333333333
$ws.Range("H5").Value = 0.61897
$ws.Range("M5").Value = 2.101131333333333
$ws.Range("N5").Value = 6.303394
$ws.Range("O5").Value = 0.03801656888058921
$ws.Range("P5").Value = 0.03801656888058921
$ws.Range("Q5").Value = 0.4335124204644445
$ws.Range("R5").Value = 3.90161178418
$ws.Range("S5").Value = 0.03801656888058921
$ws.Range("T5").Value = 0.03801656888058921
